$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column A (a row index column styled with the header/border
# style) is removed entirely; all remaining columns (B:F -> A:E) shift left.
$ws.Columns("A").Delete()
